# Sprint Report 1.docx - "Fixed a grammatical error" commit
# Applies:
#  1. Remove the _GoBack bookmark that sat before "Project Context"
#  2. (Deliverables _Toc337127350 bookmark renumbers automatically once _GoBack is removed)
#  3-10. Remove now-unneeded w:proofErr spell-check markers and merge the runs they used to
#        straddle, fix "vSTB"/" Series" -> "vSTB Series", fix "follow:" -> "follows:" grammar,
#        and re-add the _GoBack bookmark at its new location (end of the "Work that is carried"
#        heading) to match what a real Word editing session leaves behind.

$d = $word.ActiveDocument

function Set-ParaXml {
    param([string]$FindText, [string]$BodyXml)

    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $FindText
    $found = $find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $FindText"
        return
    }
    $rng = $find.Parent
    $rng.Text = ""
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $BodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# --- 1. Drop the stray _GoBack bookmark before "Project Context" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. "LoadSubscribers" no longer wrapped in proofErr spell-check markers ---
$body2 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">This was needed to add support to the </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>LoadSubscribers</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> application to mass provision these new Mobile Access features dictated by the new service GUIDs</w:t></w:r>' +
    '</w:p>'
Set-ParaXml "This was needed to add support to the LoadSubscribers application to mass provision these new Mobile Access features dictated by the new service GUIDs" $body2

# --- 3. "CurrentlyWatching" no longer wrapped in proofErr spell-check markers ---
$body3 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Implemented </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>CurrentlyWatching</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Command</w:t></w:r>' +
    '</w:p>'
Set-ParaXml "Implemented CurrentlyWatching Command" $body3

# --- 4. "Started Migration of Smartphone Emulator (Loadtest) to Emulator Controller:" one run ---
$body4 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Started Migration of Smartphone Emulator (Loadtest) to Emulator Controller:</w:t></w:r>' +
    '</w:p>'
Set-ParaXml "Started Migration of Smartphone Emulator (Loadtest) to Emulator Controller:" $body4

# --- 5. "Started converting GUI Controls to the DevExpress Equivalents:" - merge tail runs ---
$body5 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Started converting </w:t></w:r>' +
    '<w:r><w:t>GUI Controls</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>to the DevExpress Equivalents:</w:t></w:r>' +
    '</w:p>'
Set-ParaXml "Started converting GUI Controls to the DevExpress Equivalents:" $body5

# --- 6. "vSTB" + " Series" merge into one italic run "vSTB Series" ---
$body6 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Converted all of the pages listed under the </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>vSTB Series</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> group.</w:t></w:r>' +
    '</w:p>'
Set-ParaXml "Converted all of the pages listed under the vSTB Series group." $body6

# --- 7. "EmulatorController" no longer wrapped in proofErr; trailing space merges into italic run ---
$body7 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Created standalone program based off TCP code from </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">EmulatorController </w:t></w:r>' +
    '<w:r><w:t>that will be started from the Client using the WPF service</w:t></w:r>' +
    '</w:p>'
Set-ParaXml "Created standalone program based off TCP code from EmulatorController that will be started from the Client using the WPF service" $body7

# --- 8. "follow:" -> "follows:" grammar fix, split run, and the _GoBack bookmark moves here ---
$body8 = '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Work that is carried </w:t></w:r>' +
    '<w:r><w:t>over into sprint 2 is as follows:</w:t></w:r>' +
    '<w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/>' +
    '</w:p>'
Set-ParaXml "Work that is carried over into sprint 2 is as follow: " $body8

# --- 9. "Migrate Smartphone Emulator (Loadtest) to Emulator Controller" - one run ---
$rpr9 = '<w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'
$body9 = '<w:p><w:pPr><w:pStyle w:val="Default"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>' + $rpr9 + '</w:pPr>' +
    '<w:r>' + $rpr9 + '<w:t>Migrate Smartphone Emulator (Loadtest) to Emulator Controller</w:t></w:r>' +
    '</w:p>'
Set-ParaXml "Migrate Smartphone Emulator (Loadtest) to Emulator Controller" $body9

# --- 10. "Finish Converting GUI Controls to DevExpress" - one run ---
$body10 = '<w:p><w:pPr><w:pStyle w:val="Default"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>' + $rpr9 + '</w:pPr>' +
    '<w:r>' + $rpr9 + '<w:t>Finish Converting GUI Controls to DevExpress</w:t></w:r>' +
    '</w:p>'
Set-ParaXml "Finish Converting GUI Controls to DevExpress" $body10

Write-Host "done so far"
